# Updates the Price (column D) and Volume(1h) (column E) figures in the
# "cryptos" symbol-list worksheet to the latest scrape, per the commit
# "Updated symbol list on Mon Feb  6 13:43:23 UTC 2023 with GitHub Actions".
#
# The source sheet stores every data cell as literal text (inline strings,
# e.g. "327.74", "-0.29%") rather than numbers/percentages, so each write
# below forces the cell's number format to Text ("@") before assigning the
# new literal string. Without that, Excel's normal type inference would
# coerce a value like "327.81" into the number 327.81 (and "-0.43%" into
# the number -0.0043 formatted as a percentage), which would change the
# cell's stored type/value and not match the source data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet


$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "327.81"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "-0.43%"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "43.96"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "1.25%"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.563"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "-0.64%"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.08046"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "-1.98%"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.912"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "-0.20%"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "4.297"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "-3.25%"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.9450"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "0.12%"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.541"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "-11.02%"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.1167"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "-4.62%"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.1845"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "-4.32%"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.09678"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "-1.28%"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.04379"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "-1.54%"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "-0.09%"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.001286"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "0.04%"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.005865"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "-4.19%"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.495"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "-0.19%"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "9.598"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "9.09%"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.1371"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "0.03%"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.2652"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "-2.76%"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.04217"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "-4.50%"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "0.26%"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.004484"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "2.46%"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.0001261"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "2.07%"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "-0.31%"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02645"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "-6.20%"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.05496"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "-4.14%"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.007572"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "-4.32%"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.1397"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "-1.40%"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.008334"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "-16.11%"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.002010"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "-4.01%"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.008623"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "-11.39%"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00006916"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "-4.97%"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "-0.29%"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.002273"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "-0.31%"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.005449"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "60.87%"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "-0.29%"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "-0.29%"
